$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new outstation rows (Aberdeen, Dundee, Edinburgh) for rows 21-44.
$ws.Range("A21").Value = "tribunalAberdeenAddressLine1"
$ws.Range("B21").Value = "Ground Floor"
$ws.Range("B21").WrapText = $true

$ws.Range("A22").Value = "tribunalAberdeenAddressLine2"
$ws.Range("B22").Value = "AB1, 48 Huntly Street"
$ws.Range("B22").WrapText = $true

$ws.Range("A23").Value = "tribunalAberdeenTown"
$ws.Range("B23").Value = "Aberdeen"
$ws.Range("B23").WrapText = $true

$ws.Range("A24").Value = "tribunalAberdeenPostCode"
$ws.Range("B24").Value = "AB10 1SH"
$ws.Range("B24").WrapText = $true

$ws.Range("A25").Value = "tribunalAberdeenTelephone"
$ws.Range("B25").Value = "01224 593 137"

$ws.Range("A26").Value = "tribunalAberdeenFax"
$ws.Range("B26").Value = "0870 761 7766"

$ws.Range("A27").Value = "tribunalAberdeenDX"
$ws.Range("B27").Value = "DX AB77"

$ws.Range("A28").Value = "tribunalAberdeenEmail"
$ws.Range("B28").Value = "aberdeenet@justice.gov.uk"

$ws.Range("A29").Value = "tribunalDundeeAddressLine1"
$ws.Range("B29").Value = "Ground Floor"
$ws.Range("B29").WrapText = $true

$ws.Range("A30").Value = "tribunalDundeeAddressLine2"
$ws.Range("B30").Value = "Block C, Caledonian House"
$ws.Range("B30").WrapText = $true

$ws.Range("A31").Value = "tribunalDundeeAddressLine3"
$ws.Range("B31").Value = "Greenmarket"
$ws.Range("B31").WrapText = $true

$ws.Range("A32").Value = "tribunalDundeeTown"
$ws.Range("B32").Value = "Dundee"
$ws.Range("B32").WrapText = $true

$ws.Range("A33").Value = "tribunalDundeePostCode"
$ws.Range("B33").Value = "DD1 4QG"
$ws.Range("B33").WrapText = $true

$ws.Range("A34").Value = "tribunalDundeeTelephone"
$ws.Range("B34").Value = "01382 221 578"

$ws.Range("A35").Value = "tribunalDundeeFax"
$ws.Range("B35").Value = "01382 227 136"

$ws.Range("A36").Value = "tribunalDundeeDX"
$ws.Range("B36").Value = "DX DD51"

$ws.Range("A37").Value = "tribunalDundeeEmail"
$ws.Range("B37").Value = "dundeeet@justice.gov.uk"

$ws.Range("A38").Value = "tribunalEdinburghAddressLine1"
$ws.Range("B38").Value = "54-56 Melville Street"
$ws.Range("B38").WrapText = $true

$ws.Range("A39").Value = "tribunalEdinburghTown"
$ws.Range("B39").Value = "Edinburgh"
$ws.Range("B39").WrapText = $true

$ws.Range("A40").Value = "tribunalEdinburghPostCode"
$ws.Range("B40").Value = "EH3 7HF"
$ws.Range("B40").WrapText = $true

$ws.Range("A41").Value = "tribunalEdinburghTelephone"
$ws.Range("B41").Value = "0131 226 5584"

$ws.Range("A42").Value = "tribunalEdinburghFax"
$ws.Range("B42").Value = "0131 220 6847"

$ws.Range("A43").Value = "tribunalEdinburghDX"
$ws.Range("B43").Value = "DX ED147"

$ws.Range("A44").Value = "tribunalEdinburghEmail"
$ws.Range("B44").Value = "edinburghet@justice.gov.uk"

# Hyperlink the new outstation e-mail addresses, same as the existing Manchester/Glasgow rows.
$ws.Hyperlinks.Add($ws.Range("B28"), "mailto:aberdeenet@justice.gov.uk", "", "", "aberdeenet@justice.gov.uk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "mailto:dundeeet@justice.gov.uk", "", "", "dundeeet@justice.gov.uk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B44"), "mailto:edinburghet@justice.gov.uk", "", "", "edinburghet@justice.gov.uk") | Out-Null

# Match the saved view: scrolled down, with A20:B44 selected and B44 the active cell.
$ws.Range("A20:B44").Select()
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("B44").Activate()

